# "Generate Report for handoff"
#
# The handoff for 0fa1275d-8a76-4a32-8a00-1739d01a5ec1.md failed its
# transform, so a *new* source id (f749154f-d2c8-4517-bf89-a9515181c685.md)
# is generated and the status flips from "Ready for handoff" to
# "Handoff transform failed". Because no handoff ever actually went out
# for the new id, the per-language detail rows lose their
# "Latest Handoff File" link/value, their "Latest Handoff Datetime" resets
# to the zero date, and "Handoff Reason" flips from "Include" to "Ignored".

$wb = $excel.ActiveWorkbook

$oldFile = "0fa1275d-8a76-4a32-8a00-1739d01a5ec1.md"
$newFile = "f749154f-d2c8-4517-bf89-a9515181c685.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"

function Set-A2Hyperlink($ws, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$2') {
            $h.TextToDisplay = $text
        }
    }
}

# --- Overview sheet: rename the source file + new status in both columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newFile
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
Set-A2Hyperlink $overview $newFile

# --- Per-language detail sheets ---
$langSheets = @("zh-cn", "de-de")
foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = $newStatus
    Set-A2Hyperlink $ws $newFile

    # No handoff file went out anymore -> drop "Latest Handoff File" (C2)
    # entirely, hyperlink included.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$2') {
            $h.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # "Latest Handoff Datetime" resets to the zero date.
    $ws.Range("D2").Value = $zeroDate

    # Reason flips from Include to Ignored.
    $ws.Range("H2").Value = "Ignored"
}

Write-Host "Updated handoff report for $oldFile -> $newFile"
